$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.764.82'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.037.08'
$ws.Range('E3').Value = '  +2.81%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '380.58'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.92'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.548'
$ws.Range('E7').Value = '  +1.17%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.593'
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.86'
$ws.Range('E10').Value = '  +1.58%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0861'
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.511.05'
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.76'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.031.96'
$ws.Range('E16').Value = '  +2.80%  '
$ws.Range('E17').Value = '  -3.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.49'
$ws.Range('E18').Value = '  -15.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '51.754.22'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.49'
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.13'
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.34'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('E25').Value = '  -5.70%  '
$ws.Range('E26').Value = '  +2.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.60'
$ws.Range('E27').Value = '  +8.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.174'
$ws.Range('E28').Value = '  +5.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.27'
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.108'
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.30'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.07'
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '33.93'
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '50.53'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('E36').Value = '  +2.71%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.34'
$ws.Range('E38').Value = '  +5.27%  '
$ws.Range('E39').Value = '  +15.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.07'
$ws.Range('E40').Value = '  +2.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.87'
$ws.Range('E41').Value = '  +2.56%  '
$ws.Range('E42').Value = '  +3.30%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '127.73'
$ws.Range('E43').Value = '  +6.48%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.116'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.77'
$ws.Range('E45').Value = '  +5.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.74'
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('E47').Value = '  +3.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('E48').Value = '  +3.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.029.31'
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.335.54'
$ws.Range('E50').Value = '  +2.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0321'
$ws.Range('E51').Value = '  +0.81%  '
